# Fix bugs on medical records:
# Update the Email address for the first patient record (row 2) from
# "test1@gmail.com" to "newTest@gmail.com".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "newTest@gmail.com"
